# Minor fix to deliverable, Trung
# "Daniel Costaesa, Developer" -> "Daniel Gonzalez , Developer"
# (misspelled surname "Costaesa" corrected to "Gonzalez")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

$oldLine = "Daniel Costaesa, Developer"
$fullText = $tr.Text
$lineStart0 = $fullText.IndexOf($oldLine)
$lineStart = $lineStart0 + 1

$nameStart  = $lineStart        # "Daniel "
$lastStart  = $lineStart + 7    # "Costaesa"
$restStart  = $lineStart + 15   # ", Developer"

# 1) Split the trailing ", Developer" run into ", " and "Developer" while it
#    still reads ", Developer" (same-text splice keeps formatting, just adds
#    a run boundary after the comma+space).
$tr.Characters($restStart, 2).Text = ", "

# 2) Replace "Daniel Costaesa" (the first two runs) with "Daniel Gonzalez "
#    in one shot. Because the replaced range starts inside the "Daniel "
#    run, the freshly spliced text inherits that run's formatting, so the
#    misspelled surname's run (and its err="1" flag) is discarded.
$tr2 = $tf.TextRange
$tr2.Characters($nameStart, 15).Text = "Daniel Gonzalez "

# 3) Re-assert "Daniel " over its original span so the merged run above is
#    split back into "Daniel " and "Gonzalez " (both sharing the formatting
#    inherited in step 2).
$tr3 = $tf.TextRange
$tr3.Characters($nameStart, 7).Text = "Daniel "
